$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing condition string (R1 row, column C) so the AFTER timer changes from 980 to 100
$text11 = @'
(
(
( ( {
	"message": "TOKEN: abc, TOKEN: def, TOKEN: ghi, REGEX:[0-9a-z]{4}-[0-9a-z]{4}-[0-9a-z]{4}-[0-9a-z]{4}, VAR: var1, REGEX:Bar[0-9]{2}-sector[0-9]{2}, VAR: var2 "
} AND AFTER 00:00:00:100 {
	"message": "TOKEN:Node down, VAR:{var1}"
} ) OR ( {
	"message": "TOKEN:Node down, VAR:{var1}"
} AND BEFORE 00:00:00:876 {
	"message": "TOKEN:Node down, VAR:{var1}"
}) )
) AND
{
	"message": "TOKEN:Node down, VAR:{var1}"
}
)
OR
{
	"message": "TOKEN:Node down, VAR:{var1}"
}
'@

$ws.Range("C2").Value = $text11
# Writing into a wrap-text cell makes the engine recompute an autofit row
# height that doesn't match Excel's own metrics; restore the original height.
$ws.Rows.Item(2).RowHeight = 138.75

# Add new rule row 4 (R3 / Node Down) - write in the same column order the original
# author used (A, B, D, C) so new shared strings land at the same indices as the diff.
$ws.Range("A4").Value = "R3"
$ws.Range("B4").Value = "Node Down"
$solution3 = @'
Solution3: R3:
1. Node Down Errror.
'@

$ws.Range("D4").Value = $solution3
$condition3 = @'
( {
	"message": "TOKEN: abc, TOKEN: def, TOKEN: ghi, REGEX:[0-9a-z]{4}-[0-9a-z]{4}-[0-9a-z]{4}-[0-9a-z]{4}, VAR: var1, REGEX:Bar[0-9]{2}-sector[0-9]{2}, VAR: var2 "
} AND AFTER 00:00:00:100 {
	"message": "TOKEN:Node down, VAR:{var1}"
} )
'@

$ws.Range("C4").Value = $condition3

# Match formatting of the other data rows: wrap text on B:D, row height 90
$ws.Range("B4:D4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 90

# Update the view: selection moves to D4 and the window scrolls so row 4 is visible
[void]$ws.Range("D4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
